# Apply the cell-level edits captured in the target diff: updated prices
# and hourly-volume percentages for the cryptos table, plus the two pairs
# of rows (30/31 and 45/46) whose coin data got swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.520.91"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "2.371.19"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'309.29"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'104.59"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").Value = "'0.508"
$ws.Range("E7").Value = "  -5.49%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").Value = "'35.92"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "'53.47"
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").Value = "'0.0811"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "'6.98"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "2.744.12"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "'15.58"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").Value = "2.374.30"
$ws.Range("E17").Value = "  +3.21%  "
$ws.Range("D18").Value = "'0.809"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "43.495.38"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "'6.29"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("D21").Value = "'11.91"
$ws.Range("E21").Value = "  -4.87%  "
$ws.Range("D22").Value = "0.0₃0916"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "'68.32"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "'240.61"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.05"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("D26").Value = "'2.61"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "'25.77"
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'36.48"
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.20"
$ws.Range("E31").Value = "  -5.42%  "
$ws.Range("D32").Value = "'9.52"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "'160.62"
$ws.Range("D34").Value = "'5.25"
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'18.26"
$ws.Range("E37").Value = "  +5.95%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "'0.0739"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "'4.65"
$ws.Range("E40").Value = "  +9.53%  "
$ws.Range("D41").Value = "'1.93"
$ws.Range("E41").Value = "  +5.49%  "
$ws.Range("D42").Value = "'0.105"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").Value = "'0.113"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").Value = "'2.66"
$ws.Range("E44").Value = "  +15.78%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.032.27"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'19.71"
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "'3.13"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +7.13%  "
$ws.Range("D50").Value = "'58.06"
$ws.Range("E50").Value = "  +4.39%  "
$ws.Range("E51").Value = "  +0.18%  "
